$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Each entry: cell address, new value, and whether it is the numeric-looking
# "Price" column (D) that must be forced to Text so Excel does not silently
# reinterpret it as a number (which would also strip formatting like "0.380").
$updates = @(
    @{ Cell = 'D2'; Value = '66.005.36'; ForceText = $True }
    @{ Cell = 'E2'; Value = '  -0.84%  '; ForceText = $False }
    @{ Cell = 'D3'; Value = '3.521.35'; ForceText = $True }
    @{ Cell = 'E3'; Value = '  +0.41%  '; ForceText = $False }
    @{ Cell = 'E4'; Value = '  -0.15%  '; ForceText = $False }
    @{ Cell = 'D5'; Value = '576.52'; ForceText = $True }
    @{ Cell = 'E5'; Value = '  +4.69%  '; ForceText = $False }
    @{ Cell = 'D6'; Value = '178.98'; ForceText = $True }
    @{ Cell = 'E6'; Value = '  -5.84%  '; ForceText = $False }
    @{ Cell = 'D7'; Value = '0.636'; ForceText = $True }
    @{ Cell = 'E7'; Value = '  +4.78%  '; ForceText = $False }
    @{ Cell = 'E8'; Value = '  +0.06%  '; ForceText = $False }
    @{ Cell = 'E9'; Value = '  -0.11%  '; ForceText = $False }
    @{ Cell = 'E10'; Value = '  +5.60%  '; ForceText = $False }
    @{ Cell = 'D11'; Value = '55.21'; ForceText = $True }
    @{ Cell = 'E11'; Value = '  -1.17%  '; ForceText = $False }
    @{ Cell = 'E12'; Value = '  +1.92%  '; ForceText = $False }
    @{ Cell = 'D13'; Value = '9.24'; ForceText = $True }
    @{ Cell = 'E13'; Value = '  -1.91%  '; ForceText = $False }
    @{ Cell = 'D14'; Value = '4.082.52'; ForceText = $True }
    @{ Cell = 'E14'; Value = '  -0.12%  '; ForceText = $False }
    @{ Cell = 'D15'; Value = '3.515.36'; ForceText = $True }
    @{ Cell = 'E15'; Value = '  +0.10%  '; ForceText = $False }
    @{ Cell = 'E16'; Value = '  +0.31%  '; ForceText = $False }
    @{ Cell = 'D17'; Value = '18.47'; ForceText = $True }
    @{ Cell = 'E17'; Value = '  +1.03%  '; ForceText = $False }
    @{ Cell = 'D18'; Value = '12.17'; ForceText = $True }
    @{ Cell = 'E18'; Value = '  +2.66%  '; ForceText = $False }
    @{ Cell = 'D19'; Value = '65.964.19'; ForceText = $True }
    @{ Cell = 'E19'; Value = '  -0.95%  '; ForceText = $False }
    @{ Cell = 'E20'; Value = '  +1.57%  '; ForceText = $False }
    @{ Cell = 'D21'; Value = '416.06'; ForceText = $True }
    @{ Cell = 'E21'; Value = '  +2.19%  '; ForceText = $False }
    @{ Cell = 'D22'; Value = '4.20'; ForceText = $True }
    @{ Cell = 'E22'; Value = '  +6.93%  '; ForceText = $False }
    @{ Cell = 'B23'; Value = 'Litecoin'; ForceText = $False }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; ForceText = $False }
    @{ Cell = 'D23'; Value = '85.81'; ForceText = $True }
    @{ Cell = 'E23'; Value = '  +0.53%  '; ForceText = $False }
    @{ Cell = 'B24'; Value = 'Toncoin'; ForceText = $False }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $False }
    @{ Cell = 'D24'; Value = '4.27'; ForceText = $True }
    @{ Cell = 'E24'; Value = '  +1.28%  '; ForceText = $False }
    @{ Cell = 'D25'; Value = '12.86'; ForceText = $True }
    @{ Cell = 'E25'; Value = '  +7.90%  '; ForceText = $False }
    @{ Cell = 'D26'; Value = '10.95'; ForceText = $True }
    @{ Cell = 'E26'; Value = '  -1.63%  '; ForceText = $False }
    @{ Cell = 'D27'; Value = '2.86'; ForceText = $True }
    @{ Cell = 'E27'; Value = '  -2.24%  '; ForceText = $False }
    @{ Cell = 'D28'; Value = '9.03'; ForceText = $True }
    @{ Cell = 'E28'; Value = '  +1.98%  '; ForceText = $False }
    @{ Cell = 'D29'; Value = '30.45'; ForceText = $True }
    @{ Cell = 'E29'; Value = '  +0.21%  '; ForceText = $False }
    @{ Cell = 'D30'; Value = '630.73'; ForceText = $True }
    @{ Cell = 'E30'; Value = '  -4.89%  '; ForceText = $False }
    @{ Cell = 'D31'; Value = '6.43'; ForceText = $True }
    @{ Cell = 'E31'; Value = '  -4.04%  '; ForceText = $False }
    @{ Cell = 'D32'; Value = '11.68'; ForceText = $True }
    @{ Cell = 'E32'; Value = '  -0.69%  '; ForceText = $False }
    @{ Cell = 'D33'; Value = '0.111'; ForceText = $True }
    @{ Cell = 'E33'; Value = '  -0.54%  '; ForceText = $False }
    @{ Cell = 'D34'; Value = '59.81'; ForceText = $True }
    @{ Cell = 'E34'; Value = '  +0.07%  '; ForceText = $False }
    @{ Cell = 'E35'; Value = '  +11.68%  '; ForceText = $False }
    @{ Cell = 'E36'; Value = '  -0.45%  '; ForceText = $False }
    @{ Cell = 'E37'; Value = '  +0.26%  '; ForceText = $False }
    @{ Cell = 'D38'; Value = '37.38'; ForceText = $True }
    @{ Cell = 'E38'; Value = '  -3.54%  '; ForceText = $False }
    @{ Cell = 'B39'; Value = 'TheGraph'; ForceText = $False }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; ForceText = $False }
    @{ Cell = 'D39'; Value = '0.380'; ForceText = $True }
    @{ Cell = 'E39'; Value = '  -3.12%  '; ForceText = $False }
    @{ Cell = 'B40'; Value = 'Maker'; ForceText = $False }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; ForceText = $False }
    @{ Cell = 'D40'; Value = '3.274.79'; ForceText = $True }
    @{ Cell = 'E40'; Value = '  +8.97%  '; ForceText = $False }
    @{ Cell = 'D41'; Value = '3.35'; ForceText = $True }
    @{ Cell = 'E41'; Value = '  -1.07%  '; ForceText = $False }
    @{ Cell = 'D42'; Value = '0.999'; ForceText = $True }
    @{ Cell = 'E42'; Value = '  -0.40%  '; ForceText = $False }
    @{ Cell = 'E43'; Value = '  -3.83%  '; ForceText = $False }
    @{ Cell = 'E44'; Value = '  +0.92%  '; ForceText = $False }
    @{ Cell = 'B45'; Value = 'Fetch.AI'; ForceText = $False }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; ForceText = $False }
    @{ Cell = 'D45'; Value = '2.50'; ForceText = $True }
    @{ Cell = 'E45'; Value = '  -5.37%  '; ForceText = $False }
    @{ Cell = 'B46'; Value = 'ApeXProtocol'; ForceText = $False }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'; ForceText = $False }
    @{ Cell = 'D46'; Value = '3.25'; ForceText = $True }
    @{ Cell = 'E46'; Value = '  -2.27%  '; ForceText = $False }
    @{ Cell = 'E47'; Value = '  -0.24%  '; ForceText = $False }
    @{ Cell = 'E48'; Value = '  +2.25%  '; ForceText = $False }
    @{ Cell = 'B49'; Value = 'THORChain'; ForceText = $False }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; ForceText = $False }
    @{ Cell = 'D49'; Value = '8.55'; ForceText = $True }
    @{ Cell = 'E49'; Value = '  -5.42%  '; ForceText = $False }
    @{ Cell = 'B50'; Value = 'Monero'; ForceText = $False }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $False }
    @{ Cell = 'D50'; Value = '138.84'; ForceText = $True }
    @{ Cell = 'E50'; Value = '  +0.51%  '; ForceText = $False }
    @{ Cell = 'E51'; Value = '  -6.32%  '; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Preserve the existing style (border/font/alignment) while writing as text:
        # set a text number format just long enough to stop Excel from coercing the
        # literal into a number/date, then restore the original style afterwards.
        $origStyle = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = $origStyle
    } else {
        $rng.Value = $u.Value
    }
}
